$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - Court hearing
$ws.Range("C2").Value = "['court', 'hearing']"
$ws.Range("E2").Value = "['court,hearing']"

# Row 3 - Project begin
$ws.Range("C3").Value = "['not local', 'not from area', 'outside', 'only be local', 'only local people']"
$ws.Range("E3").Value = "['local,not', 'area,outside,only', 'local,only', 'local', 'people']"

# Row 4 - National Strike
$ws.Range("C4").Value = "['COSATU', 'NUMSA', 'national strike']"
$ws.Range("E4").Value = "['cosatu,numsa,national', 'strike']"

# Row 5 - Dismissals
$ws.Range("C5").Value = "['Fired', 'suspended', 'dismissed', 'discipline']"
$ws.Range("E5").Value = "['fired,suspended,dismissed,discipline']"

# Row 6 - Wage disputes
$ws.Range("C6").Value = "['Payment', 'salary', 'cheque', 'overtime', 'wage', 'wages', 'remuneration', 'not paid', 'bonus', 'bonusses', 'percent', 'unfair']"
$ws.Range("E6").Value = "['payment,salary,cheque,overtime,wage,wages,remuneration,not', 'paid,bonus,bonusses,percent,unfair']"

# Row 7 - Evictions
$ws.Range("C7").Value = "['Evicted', 'relocated', 'red ants', 'demolition', 'relocation', 'removal', 'demolished']"
$ws.Range("E7").Value = "['evicted,relocated,red', 'ants,demolition,relocation,removal,demolished']"

# Row 8 - Utility connections
$ws.Range("C8").Value = "['Disconnected', 'power is down', 'electricity is down', 'power']"
$ws.Range("D8").Value = "['cut', 'electricity', 'water', 'prepaid', 'loadshedding', 'blackout', 'watershedding']"
$ws.Range("E8").Value = "['disconnected,power', 'down,electricity', 'down,power']"
$ws.Range("F8").Value = "['cut,electricity,water,prepaid,loadshedding,blackout,watershedding']"

# Row 9 - Election related
$ws.Range("C9").Value = "['Party list', 'circulate', 'for councillor', 'the candidate', 'on the list', 'wrong candidate', 'want another person']"
$ws.Range("E9").Value = "['party', 'list,circulate,for', 'councillor,the', 'candidate,on', 'list,wrong', 'candidate,want', 'another', 'person']"

# Row 10 - Non consultation
$ws.Range("C10").Value = "['did not arrive', 'failed to respond', 'did not come back', 'did not come']"
$ws.Range("E10").Value = "['arrive,failed', 'respond,did', 'come', 'back,did', 'come']"

# Row 11 - Working conditions
$ws.Range("C11").Value = "['food', 'conditions', 'health', 'safety', 'equipment']"
$ws.Range("E11").Value = "['food,conditions,health,safety,equipment']"

# Row 12 - Arrests
$ws.Range("C12").Value = "['Arrest', 'release']"
$ws.Range("E12").Value = "['arrest,release']"

# Row 13 - Crime Event
$ws.Range("C13").Value = "['Crime', 'murder', 'kidnapped', 'was attacked', 'killed']"
$ws.Range("E13").Value = "['crime,murder,kidnapped,was', 'attacked,killed']"

# Row 14 - Neoliberal logic
$ws.Range("C14").Value = "['Capitalism', 'privatisation', 'prvatise', 'privatize', 'privatization']"
$ws.Range("E14").Value = "['capitalism,privatisation,prvatise,privatize,privatization']"
